$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.315.17'
$ws.Range("E2").Value = '  -4.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.315.82'
$ws.Range("E3").Value = '  -5.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.55'
$ws.Range("E5").Value = '  -3.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.51'
$ws.Range("E6").Value = '  -3.86%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.317.15'
$ws.Range("E8").Value = '  -5.12%  '

$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.38'
$ws.Range("E10").Value = '  -3.80%  '

$ws.Range("E11").Value = '  -4.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("E12").Value = '  -2.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.885.53'
$ws.Range("E13").Value = '  -5.10%  '

$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.316.98'
$ws.Range("E15").Value = '  -5.15%  '

$ws.Range("E16").Value = '  -5.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.86'
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.439.15'
$ws.Range("E18").Value = '  -4.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.50'
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.67'
$ws.Range("E20").Value = '  -1.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.03'
$ws.Range("E21").Value = '  -9.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '352.92'
$ws.Range("E22").Value = '  -8.64%  '

$ws.Range("E23").Value = '  -4.06%  '

$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.448.62'
$ws.Range("E25").Value = '  -5.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.55'
$ws.Range("E26").Value = '  -6.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000108'
$ws.Range("E27").Value = '  -6.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.20'
$ws.Range("E29").Value = '  -1.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.45'
$ws.Range("E30").Value = '  -3.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -3.40%  '

$ws.Range("E32").Value = '  -6.13%  '

$ws.Range("E33").Value = '  -2.88%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.351.28'
$ws.Range("E35").Value = '  -4.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.60'
$ws.Range("E36").Value = '  -2.83%  '

$ws.Range("E37").Value = '  -2.24%  '

$ws.Range("E38").Value = '  -0.69%  '

$ws.Range("E39").Value = '  -3.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '160.12'
$ws.Range("E40").Value = '  -2.64%  '

$ws.Range("E41").Value = '  -2.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("E43").Value = '  -0.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.99'
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("E45").Value = '  -7.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.11'
$ws.Range("E46").Value = '  -5.73%  '

$ws.Range("E47").Value = '  -5.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.38'
$ws.Range("E48").Value = '  -8.10%  '

$ws.Range("E49").Value = '  -1.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.869'
$ws.Range("E50").Value = '  -5.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.28'
$ws.Range("E51").Value = '  +1.76%  '

